$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999087830793
$ws.Range("A2").Value = 0.99756297670234428
$ws.Range("A3").Value = 0.98813300226291112
$ws.Range("A4").Value = 0.99087626876537505
$ws.Range("A5").Value = 0.98046471930147561
$ws.Range("A6").Value = 0.95575150573932932
$ws.Range("A7").Value = 0.94915399587762206
$ws.Range("A8").Value = 0.94002264629820209
$ws.Range("A9").Value = 0.92827556060755789
$ws.Range("A10").Value = 0.91744631229901363
$ws.Range("A11").Value = 0.91585006977033134
$ws.Range("A12").Value = 0.91306997256192268
$ws.Range("A13").Value = 0.90178194849833604
$ws.Range("A14").Value = 0.89761444596385009
$ws.Range("A15").Value = 0.89502289790214362
$ws.Range("A16").Value = 0.89251634920979572
$ws.Range("A17").Value = 0.8888083854817368
$ws.Range("A18").Value = 0.88769948401864807
$ws.Range("A19").Value = 0.99650979909830051
$ws.Range("A20").Value = 0.98939291788456663
$ws.Range("A21").Value = 0.98799445365075411
$ws.Range("A22").Value = 0.98672995185174583
$ws.Range("A23").Value = 0.98109851028778117
$ws.Range("A24").Value = 0.96807793174949419
$ws.Range("A25").Value = 0.96162099281671887
$ws.Range("A26").Value = 0.95410700332320175
$ws.Range("A27").Value = 0.95206482650867985
$ws.Range("A28").Value = 0.94392235531391888
$ws.Range("A29").Value = 0.93856912369846557
$ws.Range("A30").Value = 0.93683705276543761
$ws.Range("A31").Value = 0.94191104634934897
$ws.Range("A32").Value = 0.94432540797983866
$ws.Range("A33").Value = 0.94916724469054481
